# Update TPM-derived NATMI metrics for the Plau-Igf2r sheet.
# New ligand (Plau) average/total expression values per *sending* cluster
# and new receptor (Igf2r) average/total expression values per *target*
# cluster, together with all their derived-specificity and edge-weight
# columns, which are recomputed from those base numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand expression values, keyed by sending cluster (column A)
$ligandAvg = @{ "ECs" = 2.996089666666666; "FAPs" = 143.4723713333333; "MuSCs" = 9.885847333333333 }
$ligandTot = @{ "ECs" = 8.988268999999999; "FAPs" = 430.417114;        "MuSCs" = 29.657542 }

# New receptor expression values, keyed by target cluster (column D)
$receptorAvg = @{ "ECs" = 10.61298733333333; "FAPs" = 42.26455300000001; "MuSCs" = 24.50508366666667 }
$receptorTot = @{ "ECs" = 31.838962;         "FAPs" = 126.793659;       "MuSCs" = 73.51525100000001 }

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

# First pass: write ligand/receptor average & total expression values,
# and compute the raw (unnormalized) edge weights for every data row.
$edgeAvg = @{}
$edgeTot = @{}
$sumLigandAvg = 0
$sumLigandTot = 0
$sumEdgeAvg = 0
$sumEdgeTot = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $send = $ws.Cells.Item($r, 1).Value2
    $targ = $ws.Cells.Item($r, 4).Value2

    $g = $ligandAvg[$send]
    $h = $ligandTot[$send]
    $m = $receptorAvg[$targ]
    $n = $receptorTot[$targ]

    $ws.Cells.Item($r, 7).Value2  = $g   # G: Ligand average expression value
    $ws.Cells.Item($r, 8).Value2  = $h   # H: Ligand total expression value
    $ws.Cells.Item($r, 13).Value2 = $m   # M: Receptor average expression value
    $ws.Cells.Item($r, 14).Value2 = $n   # N: Receptor total expression value

    $q = $g * $m
    $t = $h * $n
    $edgeAvg[$r] = $q
    $edgeTot[$r] = $t
    $sumEdgeAvg += $q
    $sumEdgeTot += $t
}

# Sum of ligand average/total expression values across all sending clusters
# (each cluster counted once) is needed to normalize ligand specificity.
$uniqueSend = @{}
$uniqueTarg = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $send = $ws.Cells.Item($r, 1).Value2
    $targ = $ws.Cells.Item($r, 4).Value2
    if (-not $uniqueSend.ContainsKey($send)) { $uniqueSend[$send] = $true }
    if (-not $uniqueTarg.ContainsKey($targ)) { $uniqueTarg[$targ] = $true }
}
foreach ($k in $uniqueSend.Keys) {
    $sumLigandAvg += $ligandAvg[$k]
    $sumLigandTot += $ligandTot[$k]
}
$sumReceptorAvg = 0
$sumReceptorTot = 0
foreach ($k in $uniqueTarg.Keys) {
    $sumReceptorAvg += $receptorAvg[$k]
    $sumReceptorTot += $receptorTot[$k]
}

# Second pass: derived specificities (ligand I/J, receptor O/P) and
# edge-weight derived specificities (S/T), normalized across all rows.
for ($r = 2; $r -le $lastRow; $r++) {
    $send = $ws.Cells.Item($r, 1).Value2
    $targ = $ws.Cells.Item($r, 4).Value2

    $i = $ligandAvg[$send] / $sumLigandAvg
    $j = $ligandTot[$send] / $sumLigandTot
    $o = $receptorAvg[$targ] / $sumReceptorAvg
    $p = $receptorTot[$targ] / $sumReceptorTot

    $ws.Cells.Item($r, 9).Value2  = $i   # I: Ligand derived specificity (avg)
    $ws.Cells.Item($r, 10).Value2 = $j   # J: Ligand derived specificity (total)
    $ws.Cells.Item($r, 15).Value2 = $o   # O: Receptor derived specificity (avg)
    $ws.Cells.Item($r, 16).Value2 = $p   # P: Receptor derived specificity (total)

    $q = $edgeAvg[$r]
    $t = $edgeTot[$r]
    $s = $q / $sumEdgeAvg
    $u = $t / $sumEdgeTot

    $ws.Cells.Item($r, 17).Value2 = $q   # Q: Edge average expression weight
    $ws.Cells.Item($r, 18).Value2 = $t   # R: Edge total expression weight
    $ws.Cells.Item($r, 19).Value2 = $s   # S: Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value2 = $u   # T: Edge total expression derived specificity
}
